$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text in columns D/E is written back as literal text
# (matches the source data, which stores prices/deltas as inline strings).

# Row 2: 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.888.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.59%  "

# Row 3: 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.864.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.60%  "

# Row 4: 'TetherUSD'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "

# Row 5: 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.31%  "

# Row 6: 'Solana'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.27"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.73%  "

# Row 7: 'XRP'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.671"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.46%  "

# Row 9: 'Cardano'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.752"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.60%  "

# Row 10: 'Dogecoin'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.175"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.18%  "

# Row 11: 'Avalanche'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.44"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.09%  "

# Row 12: 'ShibaInu'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000320"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.81%  "

# Row 13: 'Polkadot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.46"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.96%  "

# Row 14: 'WrappedliquidstakedEther2.0'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.487.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.56%  "

# Row 15: 'Chainlink' -> 'WrappedEther'
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.887.04"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.90%  "

# Row 16: 'WrappedEther' -> 'Chainlink'
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.67%  "

# Row 17: 'Uniswap'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.77"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.91%  "

# Row 18: 'Polygon'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.62%  "

# Row 19: 'TRON'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.13%  "

# Row 20: 'WrappedBTC'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.921.76"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.15%  "

# Row 21: 'BitcoinCash'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "436.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.63%  "

# Row 22: 'PancakeSwap'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.71"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.99%  "

# Row 23: 'Litecoin'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.33"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.72%  "

# Row 24: 'ImmutableX'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.27"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.58%  "

# Row 25: 'InternetComputer(DFINITY)'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.87"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.35%  "

# Row 26: 'Toncoin' -> 'RenderToken'
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.34"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.73%  "

# Row 27: 'RenderToken' -> 'Toncoin'
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.41%  "

# Row 28: 'LEO'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.92"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.22%  "

# Row 29: 'Filecoin'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.35"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.75%  "

# Row 30: 'EthereumClassic'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.10"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.47%  "

# Row 31: 'NEARProtocol'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.43%  "

# Row 32: 'Cosmos'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.57"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.73%  "

# Row 33: 'InjectiveProtocol'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "48.67"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.30%  "

# Row 34: 'Hedera'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.10%  "

# Row 35: 'PEPE' -> 'OKB'
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "69.40"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.20%  "

# Row 36: 'OKB' -> 'PEPE'
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0987"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +13.69%  "

# Row 37: 'Bittensor'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "632.83"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.85%  "

# Row 38: 'TheGraph'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.429"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.09%  "

# Row 39: 'Kaspa'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.75%  "

# Row 40: 'Dai'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41: 'ThetaToken' -> 'FirstDigitalUSD'
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.06%  "

# Row 42: 'dogwifhat' -> 'ThetaToken'
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.28"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.67%  "

# Row 43: 'FirstDigitalUSD' -> 'dogwifhat'
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +25.78%  "

# Row 44: 'VeChain'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0470"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.35%  "

# Row 45: 'THORChain'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.10"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.54%  "

# Row 46: 'Fetch.AI'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.34%  "

# Row 47: 'Stellar'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.56%  "

# Row 48: 'WEMIXToken' -> 'ApeXProtocol'
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.22%  "

# Row 49: 'ApeXProtocol' -> 'WEMIXToken'
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -15.14%  "

# Row 50: 'Maker'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.839.49"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.27%  "

# Row 51: 'FLOKI'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000272"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.13%  "
